$wb = $excel.ActiveWorkbook

# Update status text from "Ready for handoff" to "In Translation"
$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("E2:F3").Value = "In Translation"
$zhcn.Range("C2:C3").Value = "In Translation"
$dede.Range("C2:C3").Value = "In Translation"

# Narrow the "Status" columns (previously widened to fit "Ready for handoff").
# Target stored width is 13.4101848602295 characters; the engine quantizes
# ColumnWidth assignments to 1/6-character steps, so 12.5 is the closest
# achievable input (-> stored width 13.3333..., nearest representable value).
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
